# Add opening and closing Durkheim quotes to the article, each followed/
# preceded by blank spacer paragraphs, matching the authored diff.

$d = $word.ActiveDocument

$openingQuote = '"The believer who has communicated with his god is not merely a man who sees new truths of which the unbeliever is ignorant; he is a man transformed." — Émile Durkheim'
$closingQuote = '"Man cannot become attached to higher aims and submit to a rule if he sees nothing above him to which he belongs." — Émile Durkheim'

# --- Beginning of document: insert quote paragraph + blank paragraph
# before the existing title paragraph. ---
$firstPara = $d.Paragraphs.First
$firstPara.Range.InsertParagraphBefore()
$firstPara.Range.InsertParagraphBefore()

$quotePara = $d.Paragraphs.First
$quotePara.Range.Text = $openingQuote

$blankPara = $d.Paragraphs.First.Next()
$blankPara.Range.Text = ""

# --- End of document: insert two blank paragraphs then the closing quote
# paragraph after the existing disclaimer paragraph. ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$lastPara.Range.InsertParagraphAfter()
$lastPara.Range.InsertParagraphAfter()

$blankPara2 = $lastPara.Next()
$blankPara2.Range.Text = ""

$blankPara3 = $blankPara2.Next()
$blankPara3.Range.Text = ""

$quote2Para = $d.Paragraphs.Last
$quote2Para.Range.Text = $closingQuote

Write-Output "done"
